# Add a "metric_collectors" field/column to the Benchmark sheet, directly
# before the existing "stages" column (I), shifting stages/id/name/description
# one column to the right (I->J, J->K, K->L, L->M).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Benchmark")

# Insert a new blank column at I, pushing "stages" (and the columns after it)
# one slot to the right.
$ws.Range("I1").EntireColumn.Insert()

# Populate the newly inserted header cell.
$ws.Range("I1").Value = "metric_collectors"
